$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings with European-style "." separators that Excel
# would otherwise auto-convert to numbers. Force the range to Text format,
# write the literal strings, then restore the default "Normal" style so the
# cells don't carry a lingering number-format override.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.498.63'
$ws.Range('D3').Value = '1.836.16'
$ws.Range('D4').Value = '1.002'
$ws.Range('D5').Value = '256.80'
$ws.Range('D7').Value = '0.5238'
$ws.Range('D8').Value = '0.3148'
$ws.Range('D9').Value = '0.06780'
$ws.Range('D10').Value = '18.66'
$ws.Range('D11').Value = '0.7743'
$ws.Range('D12').Value = '0.07766'
$ws.Range('D13').Value = '1.823.03'
$ws.Range('D14').Value = '87.56'
$ws.Range('D15').Value = '4.993'
$ws.Range('D16').Value = '1.002'
$ws.Range('D17').Value = '13.79'
$ws.Range('D19').Value = '0.000007904'
$ws.Range('D20').Value = '26.531.84'
$ws.Range('D21').Value = '2.070.23'
$ws.Range('D22').Value = '4.586'
$ws.Range('D23').Value = '5.947'
$ws.Range('D24').Value = '9.275'
$ws.Range('D26').Value = '2.203'
$ws.Range('D27').Value = '1.671'
$ws.Range('D28').Value = '16.89'
$ws.Range('D29').Value = '111.47'
$ws.Range('D30').Value = '4.147'
$ws.Range('D31').Value = '0.08721'
$ws.Range('D32').Value = '4.053'
$ws.Range('D33').Value = '0.04857'
$ws.Range('D34').Value = '1.131'
$ws.Range('D35').Value = '0.7169'
$ws.Range('D36').Value = '2.859'
$ws.Range('D37').Value = '3.086'
$ws.Range('D38').Value = '2.221'
$ws.Range('D39').Value = '0.01728'
$ws.Range('D40').Value = '0.4797'
$ws.Range('D42').Value = '109.82'
$ws.Range('D43').Value = '5.924'
$ws.Range('D44').Value = '1.001'
$ws.Range('D45').Value = '7.607'
$ws.Range('D46').Value = '0.4149'
$ws.Range('D47').Value = '8.997'
$ws.Range('D48').Value = '0.1225'
$ws.Range('D49').Value = '0.05803'
$ws.Range('D50').Value = '34.65'
$ws.Range('D51').Value = '0.8904'

$dRange.Style = "Normal"

# Column E holds percentage strings (with surrounding spaces); Excel keeps
# these as text automatically so they can be set directly.
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  -4.04%  '
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('E27').Value = '  +1.47%  '
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('E32').Value = '  -2.18%  '
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('E35').Value = '  +1.28%  '
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('E39').Value = '  -1.87%  '
$ws.Range('E40').Value = '  -1.17%  '
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('E43').Value = '  -2.59%  '
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('E45').Value = '  -1.58%  '
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('E51').Value = '  +0.30%  '
